$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = "1608894"
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 150
$ws.Range("J6").Value = 20

# Row 7
$ws.Range("A7").Value = "106241"
$ws.Range("F7").Value = 60
$ws.Range("G7").Value = 105
$ws.Range("H7").Value = 40
$ws.Range("I7").Value = 220
$ws.Range("J7").Value = 75

# Row 8
$ws.Range("A8").Value = "1608883"
$ws.Range("F8").Value = 245
$ws.Range("G8").Value = 20
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 90
$ws.Range("J8").Value = 155

# Row 9
$ws.Range("A9").Value = "63366"
$ws.Range("D9").Value = "Pickup"
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 465
$ws.Range("J9").Value = 0

# Row 10
$ws.Range("A10").Value = "10201"
$ws.Range("D10").Value = "Pickup"
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 485
$ws.Range("J10").Value = 100

# Row 11
$ws.Range("A11").Value = "101418"
$ws.Range("D11").Value = "Pickup"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 90
$ws.Range("H11").Value = 150
$ws.Range("I11").Value = 365
$ws.Range("J11").Value = 0

# Row 12
$ws.Range("A12").Value = "15025"
$ws.Range("D12").Value = "Pickup"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 110
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 0

# Row 13
$ws.Range("A13").Value = "47012"
$ws.Range("D13").Value = "Pickup"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 450
$ws.Range("J13").Value = 0

# Row 14
$ws.Range("A14").Value = "15025"
$ws.Range("D14").Value = "Pickup"
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 110
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 0

# Row 15
$ws.Range("A15").Value = "47012"
$ws.Range("D15").Value = "Pickup"
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 450
$ws.Range("J15").Value = 0

# Row 16
$ws.Range("A16").Value = "34000"
$ws.Range("D16").Value = "Pickup"
$ws.Range("F16").Value = 285
$ws.Range("G16").Value = 110
$ws.Range("H16").Value = 30
$ws.Range("I16").Value = 510
$ws.Range("J16").Value = 0

# Row 17
$ws.Range("A17").Value = "226547"
$ws.Range("D17").Value = "Pickup"
$ws.Range("F17").Value = 1100
$ws.Range("G17").Value = 130
$ws.Range("H17").Value = 700
$ws.Range("I17").Value = 1500
$ws.Range("J17").Value = 0

# Row 18
$ws.Range("A18").Value = "1609005"
$ws.Range("D18").Value = "Pickup"
$ws.Range("F18").Value = 500
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 400
$ws.Range("J18").Value = 0

# Row 19
$ws.Range("A19").Value = "1001188"
$ws.Range("D19").Value = "Pickup"
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 285
$ws.Range("J19").Value = 0

# Row 20
$ws.Range("A20").Value = "1001265"
$ws.Range("D20").Value = "Pickup"
$ws.Range("F20").Value = 80
$ws.Range("G20").Value = 150
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 290
$ws.Range("J20").Value = 0

# Row 21
$ws.Range("A21").Value = "291066"
$ws.Range("D21").Value = "Pickup"
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 20
$ws.Range("J21").Value = 305

# Row 22
$ws.Range("A22").Value = "100234"
$ws.Range("D22").Value = "Pickup"
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 855
$ws.Range("J22").Value = 35

# Row 23
$ws.Range("A23").Value = "1608719"
$ws.Range("D23").Value = "Pickup"
$ws.Range("F23").Value = 205
$ws.Range("G23").Value = 135
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 195
$ws.Range("J23").Value = 515

# Row 24
$ws.Range("A24").Value = "147832"
$ws.Range("D24").Value = "Pickup"
$ws.Range("F24").Value = 2735
$ws.Range("G24").Value = 225
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 2120
$ws.Range("J24").Value = 0

# Row 25
$ws.Range("A25").Value = "1608526"
$ws.Range("D25").Value = "Pickup"
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 365
$ws.Range("H25").Value = 95
$ws.Range("I25").Value = 220
$ws.Range("J25").Value = 0

# Row 26
$ws.Range("A26").Value = "1005313"
$ws.Range("D26").Value = "Pickup"
$ws.Range("F26").Value = 800
$ws.Range("G26").Value = 150
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 250
$ws.Range("J26").Value = 410

# Row 27
$ws.Range("A27").Value = "102616"
$ws.Range("D27").Value = "Pickup"
$ws.Range("F27").Value = 200
$ws.Range("G27").Value = 75
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 150
$ws.Range("J27").Value = 0

# Row 28
$ws.Range("A28").Value = "100243"
$ws.Range("D28").Value = "Pickup"
$ws.Range("F28").Value = 900
$ws.Range("G28").Value = 200
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 400
$ws.Range("J28").Value = 0

# Row 29
$ws.Range("A29").Value = "23004"
$ws.Range("D29").Value = "Pickup"
$ws.Range("F29").Value = 75
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0

# Row 30
$ws.Range("A30").Value = "18383"
$ws.Range("D30").Value = "Pickup"
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 350

# Row 31
$ws.Range("A31").Value = "6338555"
$ws.Range("D31").Value = "Pickup"
$ws.Range("F31").Value = 900
$ws.Range("G31").Value = 200
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 400
$ws.Range("J31").Value = 0

# Row 67
$ws.Range("A67").Value = "1147"
$ws.Range("F67").Value = 305
$ws.Range("G67").Value = 125
$ws.Range("H67").Value = 50
$ws.Range("I67").Value = 925
$ws.Range("J67").Value = 250

# Row 68
$ws.Range("A68").Value = "364"
$ws.Range("D68").Value = "Deliver"
$ws.Range("F68").Value = 1885
$ws.Range("G68").Value = 550
$ws.Range("H68").Value = 880
$ws.Range("I68").Value = 6160
$ws.Range("J68").Value = 100

# Row 69
$ws.Range("A69").Value = "483"
$ws.Range("D69").Value = "Deliver"
$ws.Range("F69").Value = 3060
$ws.Range("G69").Value = 875
$ws.Range("H69").Value = 95
$ws.Range("I69").Value = 3985
$ws.Range("J69").Value = 855

# Row 70
$ws.Range("A70").Value = "1255"
$ws.Range("D70").Value = "Deliver"
$ws.Range("F70").Value = 2875
$ws.Range("G70").Value = 625
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 1200
$ws.Range("J70").Value = 760
